$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1
$BR = [char]11

# ---------------------------------------------------------------------------
# 1) Simple whole-paragraph text swaps (single-run paragraphs, stable index)
# ---------------------------------------------------------------------------

# Objetivos (English, italic): long PBL blurb -> short "Identify..." blurb
$d.Paragraphs(7).Range.Text = "Identify a real or potential problem in a production system and propose a solution to the problem."

# Programa resumido (Portuguese): "A recuperação é contínua..." -> "Livros e Artigos..."
$d.Paragraphs(11).Range.Text = "Livros e Artigos científicos relacionados com o tema do projeto/problema."

# Programa resumido (English, italic): short "Identify..." blurb -> long PBL blurb
$d.Paragraphs(12).Range.Text = "Develop an interdisciplinary project, of medium complexity, on a topic related to Production Engineering, similar to situations that students will encounter in real life, in the actual exercise of their profession; Apply and integrate knowledge acquired in other course subjects; Develop technical skills (related to the project itself), and transversal skills (active learning, systemic thinking, problem-solving skills, teamwork, leadership, interpersonal relationships, conflict management, communication skills, planning skills, creativity and initiative), in a learning environment based on PBL (Project-Based Learning and Problem -Baed Learning)."

# Programa: "Livros e Artigos..." -> "O grupo social alvo..."
$d.Paragraphs(14).Range.Text = "O grupo social alvo são médias e grandes empresas, incluindo os profissionais dessas empresas, da Região do Vale do Paraíba."

# Bibliografia: "Será realizada uma pesquisa..." -> "11079086 - Herlandí de Souza Andrade"
$d.Paragraphs(19).Range.Text = "11079086 - Herlandí de Souza Andrade"

# ---------------------------------------------------------------------------
# 2) "Docente(s) Responsável(eis)" list paragraph restructure (Paragraphs(9))
# ---------------------------------------------------------------------------

# Remove the "11079086 - Herlandí de Souza Andrade" run (+ its trailing break)
$f = $d.Content
$f.Find.Execute("11079086 - Herlandí de Souza Andrade", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$delRange = $d.Range($f.Start, $f.End + 1)
$delRange.Delete()

# Remove the (now orphaned) "Identificar um problema..." run (+ its trailing break)
$f = $d.Content
$f.Find.Execute("Identificar um problema real ou potencial em um sistema produtivo e propor uma solução para o problema.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$delRange = $d.Range($f.Start, $f.End + 1)
$delRange.Delete()

# Re-insert "Identificar um problema..." as a new first run of the paragraph
$d.Paragraphs(9).Range.InsertBefore("Identificar um problema real ou potencial em um sistema produtivo e propor uma solução para o problema." + $BR)

# Give the "A nota será..." run a trailing break, then append the
# "A recuperação..." text as a new trailing run of the same paragraph
$f = $d.Content
$f.Find.Execute("A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$f.Text = $f.Text + $BR
$d.Paragraphs(9).Range.InsertAfter("A recuperação é contínua ao longo da disciplina, considerando as diversas atividades e entregas a serem realizadas. Não há prova de recuperação.")

# ---------------------------------------------------------------------------
# 3) "Avaliação" list paragraph restructure (Paragraphs(17))
# ---------------------------------------------------------------------------

# Remove the "O grupo social alvo..." run (+ its trailing break)
$f = $d.Content
$f.Find.Execute("O grupo social alvo são médias e grandes empresas, incluindo os profissionais dessas empresas, da Região do Vale do Paraíba.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$delRange = $d.Range($f.Start, $f.End + 1)
$delRange.Delete()

# Remove the original bold "Critério: " run entirely
$f = $d.Content
$f.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$f.Delete()

# Relabel the bold "Norma de recuperação: " run to "Critério: "
$f = $d.Content
$f.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Critério: ", $wdReplaceOne) | Out-Null

# Give the last bullet of the "Critério" body text a trailing break, then
# append a new bold "Norma de recuperação: " run and its body-text run
$f = $d.Content
$f.Find.Execute("8.Realização da avaliação do projeto pela empresa, autoavaliação pelos estudantes e lições aprendidas.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$f.Text = $f.Text + $BR

$newBold = $d.Paragraphs(17).Range
$newBold.InsertAfter("Norma de recuperação: ")
$newBold2 = $d.Content
$newBold2.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
# the Find above matches the first "Norma de recuperação: " (the relabel target no
# longer has this text), i.e. the run we just inserted at the end of Paragraphs(17)
$newBold2.Bold = 1

$d.Paragraphs(17).Range.InsertAfter("Será realizada uma pesquisa de satisfação com os responsáveis pela empresa participante da atividade, durante e após o projeto. Após a pesquisa, o grupo de estudantes da disciplina, fará uma análise dos resultados e uma autoavaliação e discutirá tais resultados com o professor da disciplina, de maneira e retroalimentar a realização do projeto.")

Write-Output "done"
